$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = 2025092219
$ws.Range("O16").Value = "服務"
$ws.Range("S16").Value = "O"
$ws.Range("B17").Value = 2025092220
$ws.Range("O17").Value = "抄表"
$ws.Range("S17").Value = ""
$ws.Range("B98").Value = 2025091374
$ws.Range("O98").Value = "抄表"
$ws.Range("Q98").Value = 29573
$ws.Range("T98").Value = ""
$ws.Range("B99").Value = 2025091372
$ws.Range("Q99").Value = 29570
$ws.Range("S99").Value = ""
$ws.Range("T99").Value = "O"
$ws.Range("AA99").Value = "O"
$ws.Range("B100").Value = 2025091373
$ws.Range("O100").Value = "服務"
$ws.Range("S100").Value = "O"
$ws.Range("AA100").Value = ""
$ws.Range("B105").Value = 2025092237
$ws.Range("O105").Value = "服務"
$ws.Range("S105").Value = "O"
$ws.Range("AA105").Value = ""
$ws.Range("AB105").Value = ""
$ws.Range("B106").Value = 2025092238
$ws.Range("O106").Value = "抄表"
$ws.Range("S106").Value = ""
$ws.Range("AA106").Value = "O"
$ws.Range("B107").Value = 2025090954
$ws.Range("O107").Value = "抄表"
$ws.Range("S107").Value = ""
$ws.Range("B108").Value = 2025090953
$ws.Range("O108").Value = "服務"
$ws.Range("S108").Value = "O"
$ws.Range("B124").Value = 2025091708
$ws.Range("O124").Value = "抄表"
$ws.Range("S124").Value = ""
$ws.Range("AA124").Value = "O"
$ws.Range("AB124").Value = 1
$ws.Range("B125").Value = 2025091707
$ws.Range("O125").Value = "服務"
$ws.Range("S125").Value = "O"
$ws.Range("AA125").Value = ""
$ws.Range("B134").Value = 2025090972
$ws.Range("O134").Value = "抄表"
$ws.Range("S134").Value = ""
$ws.Range("B135").Value = 2025090971
$ws.Range("O135").Value = "服務"
$ws.Range("S135").Value = "O"
$ws.Range("B151").Value = 2025092229
$ws.Range("O151").Value = "服務"
$ws.Range("S151").Value = "O"
$ws.Range("AA151").Value = ""
$ws.Range("AB151").Value = ""
$ws.Range("B152").Value = 2025092230
$ws.Range("O152").Value = "抄表"
$ws.Range("S152").Value = ""
$ws.Range("AA152").Value = "O"
$ws.Range("B161").Value = 2025090403
$ws.Range("O161").Value = "抄表"
$ws.Range("S161").Value = ""
$ws.Range("AA161").Value = "O"
$ws.Range("AB161").Value = 1
$ws.Range("B162").Value = 2025090402
$ws.Range("O162").Value = "服務"
$ws.Range("S162").Value = "O"
$ws.Range("AA162").Value = ""
$ws.Range("B260").Value = 2025090821
$ws.Range("J260").Value = ""
$ws.Range("K260").Value = ""
$ws.Range("O260").Value = "服務"
$ws.Range("S260").Value = "O"
$ws.Range("Y260").Value = "O"
$ws.Range("Z260").Value = "PMQ3"
$ws.Range("AA260").Value = ""
$ws.Range("AB260").Value = ""
$ws.Range("B261").Value = 2025090712
$ws.Range("J261").Value = "一般件"
$ws.Range("K261").Value = "其他"
$ws.Range("O261").Value = "維修"
$ws.Range("S261").Value = ""
$ws.Range("Y261").Value = ""
$ws.Range("Z261").Value = "更換發票機`n換上 8155006323`n換下 8155004152"
$ws.Range("AA261").Value = "O"
$ws.Range("B270").Value = 2025091612
$ws.Range("H270").Value = "14:50:00"
$ws.Range("Z270").Value = "更換發票機`n換上 8155004438`n換下 8155006284"
$ws.Range("B271").Value = 2025091812
$ws.Range("H271").Value = "15:25:00"
$ws.Range("Z271").Value = "到場後詢問門市表示機器已無問題，`n並且報修櫃號為TM2非TM1`n現場觀察半個小時以上觸控均無異狀`n1 嘗試交易 正常`n2 登出登入 正常`n3 重新開機 正常"
$ws.Range("AB278").Value = 154
